# Helper: replace a literal XML fragment inside a paragraph's WordOpenXML
# representation, then write it back via Range.InsertXML. This lets us
# precisely control run-splitting / proofErr / bookmark placement the way
# Word's own editors (spell-check, typing) would produce, without having
# to fight the higher-level Range.Text / Find.Execute API for structural
# changes.
function Set-ParaXml {
    param($para, [string]$oldFrag, [string]$newFrag)
    $r = $para.Range
    $xml = $r.WordOpenXML
    if ($xml.IndexOf($oldFrag) -lt 0) {
        throw "fragment not found in paragraph: $oldFrag"
    }
    $xml = $xml.Replace($oldFrag, $newFrag)
    $r.InsertXML($xml)
}

$d = $word.ActiveDocument

# 1) "Date: 4-18-" + "19" (two runs) -> single run "Date: 4-18-19"
$d.Content.Find.Execute("Date: 4-18-19", $true, $false, $false, $false, $false, $true, 1, $false, "Date: 4-18-19", 2) | Out-Null

# 2) "Instructor: Professor VanderLinden" -> split off "VanderLinden" with
#    spell-check proofErr markers around it (simulates Word's proofer
#    flagging the surname as a misspelling).
$p = $d.Paragraphs(4)
Set-ParaXml $p `
    '<w:r><w:t>Instructor: Professor VanderLinden</w:t></w:r>' `
    '<w:r><w:t xml:space="preserve">Instructor: Professor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VanderLinden</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# 3) "Use Keras (TensorFlow) deep neural networks..." ->
#    "Use Keras to implement deep neural networks..." with spell-check
#    proofErr markers around "Keras".
$p = $d.Paragraphs(8)
Set-ParaXml $p `
    '<w:r><w:t>Use Keras (TensorFlow) deep neural networks to do SLO topic classification over the standard TBL topics using Tweets relating to mining companies:</w:t></w:r>' `
    '<w:r><w:t xml:space="preserve">Use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Keras</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">to implement </w:t></w:r><w:r><w:t>deep neural networks to do SLO topic classification over the standard TBL topics using Tweets relating to mining companies:</w:t></w:r>'

# 4) Insert a new bold "Project Ideas:" paragraph right before
#    "If necessary, manually hand-tag..." (currently paragraph 19).
$p19 = $d.Paragraphs(19)
$p19.Range.InsertParagraphBefore()
$projectIdeas = $d.Paragraphs(19)
Set-ParaXml $projectIdeas `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r></w:r>' `
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Project Ideas:</w:t></w:r>'

# 5) Insert two new bullet paragraphs right after
#    "If necessary, manually hand-tag..." (now paragraph 20).
$ifNecessary = $d.Paragraphs(20)
$ifNecessary.Range.InsertParagraphAfter()
$gpu1 = $d.Paragraphs(21)
Set-ParaXml $gpu1 `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r></w:r>' `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Enable GPU support for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Keras</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tensorflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> back-end for training the model.</w:t></w:r>'

$gpu1Fresh = $d.Paragraphs(21)
$gpu1Fresh.Range.InsertParagraphAfter()
$gpu2 = $d.Paragraphs(22)
Set-ParaXml $gpu2 `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r></w:r>' `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Parallelize GPU support for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Keras</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tensorflow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> back-end for training the model.</w:t></w:r>'

Write-Output "done"
